$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:A250").Value = 3215996243
$ws.Range("A251:A500").Value = 3104023154
$ws.Range("A501:A750").Value = 3174466432

[void]$ws.Range("C747").Select()
$excel.ActiveWindow.ScrollRow = 735
$excel.ActiveWindow.ScrollColumn = 1
